$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 3
$ws.Range("B2").Value2 = 6576985
$ws.Range("E2").Value2 = "Cercle Brugge"
$ws.Range("F2").Value2 = "Westerlo"
$ws.Range("G2").Value2 = 2
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("L2").Value2 = 1.727
$ws.Range("M2").Value2 = 4.2
$ws.Range("N2").Value2 = 4.2
$ws.Range("O2").Value2 = 1.45
$ws.Range("P2").Value2 = 5
$ws.Range("Q2").Value2 = 5.75
$ws.Range("R2").Value2 = -1.25
$ws.Range("S2").Value2 = 1.95
$ws.Range("T2").Value2 = 1.9
$ws.Range("V2").Value2 = 1.9
$ws.Range("W2").Value2 = 1.95
$ws.Range("X2").Value2 = 0.45
$ws.Range("AA2").Value2 = 0.95
$ws.Range("AC2").Value2 = -1
$ws.Range("AD2").Value2 = 0.95

# Row 3 <- original row 2
$ws.Range("B3").Value2 = 6576986
$ws.Range("E3").Value2 = "Gent"
$ws.Range("F3").Value2 = "Standard Liege"
$ws.Range("G3").Value2 = 3
$ws.Range("H3").Value2 = 1
$ws.Range("I3").Value2 = 2
$ws.Range("J3").Value2 = 1
$ws.Range("L3").Value2 = 1.5
$ws.Range("M3").Value2 = 4.75
$ws.Range("N3").Value2 = 5.5
$ws.Range("O3").Value2 = 1.363
$ws.Range("P3").Value2 = 5.75
$ws.Range("Q3").Value2 = 6.5
$ws.Range("R3").Value2 = -1.5
$ws.Range("S3").Value2 = 1.925
$ws.Range("T3").Value2 = 1.925
$ws.Range("V3").Value2 = 1.925
$ws.Range("W3").Value2 = 1.925
$ws.Range("X3").Value2 = 0.363
$ws.Range("AA3").Value2 = 0.925
$ws.Range("AC3").Value2 = 0.925
$ws.Range("AD3").Value2 = -1

# Row 34 <- original row 35
$ws.Range("B34").Value2 = 6810015
$ws.Range("E34").Value2 = "Gent"
$ws.Range("F34").Value2 = "SintTruidense"
$ws.Range("G34").Value2 = 2
$ws.Range("H34").Value2 = 2
$ws.Range("I34").Value2 = 1
$ws.Range("J34").Value2 = 2
$ws.Range("L34").Value2 = 1.45
$ws.Range("M34").Value2 = 4.5
$ws.Range("N34").Value2 = 7
$ws.Range("O34").Value2 = 1.533
$ws.Range("P34").Value2 = 4.2
$ws.Range("Q34").Value2 = 6
$ws.Range("S34").Value2 = 1.9
$ws.Range("T34").Value2 = 1.95
$ws.Range("U34").Value2 = 2.5
$ws.Range("V34").Value2 = 1.8
$ws.Range("W34").Value2 = 2.05
$ws.Range("Y34").Value2 = 3.2
$ws.Range("AB34").Value2 = 0.95
$ws.Range("AC34").Value2 = 0.8
$ws.Range("AD34").Value2 = -1

# Row 35 <- original row 34
$ws.Range("B35").Value2 = 6810012
$ws.Range("E35").Value2 = "Genk"
$ws.Range("F35").Value2 = "Charleroi"
$ws.Range("G35").Value2 = 0
$ws.Range("H35").Value2 = 0
$ws.Range("I35").Value2 = 0
$ws.Range("J35").Value2 = 0
$ws.Range("L35").Value2 = 1.533
$ws.Range("M35").Value2 = 4.2
$ws.Range("N35").Value2 = 5.5
$ws.Range("O35").Value2 = 1.6
$ws.Range("P35").Value2 = 4
$ws.Range("Q35").Value2 = 5.25
$ws.Range("S35").Value2 = 1.975
$ws.Range("T35").Value2 = 1.875
$ws.Range("U35").Value2 = 3
$ws.Range("V35").Value2 = 1.9
$ws.Range("W35").Value2 = 1.95
$ws.Range("Y35").Value2 = 3
$ws.Range("AB35").Value2 = 0.875
$ws.Range("AC35").Value2 = -1
$ws.Range("AD35").Value2 = 0.95

# Row 155 <- original row 156
$ws.Range("B155").Value2 = 6810130
$ws.Range("E155").Value2 = "Antwerp"
$ws.Range("F155").Value2 = "Westerlo"
$ws.Range("G155").Value2 = 2
$ws.Range("H155").Value2 = 2
$ws.Range("I155").Value2 = 1
$ws.Range("K155").Value2 = "D"
$ws.Range("L155").Value2 = 1.363
$ws.Range("M155").Value2 = 5
$ws.Range("N155").Value2 = 7.5
$ws.Range("O155").Value2 = 1.333
$ws.Range("P155").Value2 = 5.25
$ws.Range("Q155").Value2 = 8
$ws.Range("R155").Value2 = -1.5
$ws.Range("S155").Value2 = 1.9
$ws.Range("T155").Value2 = 1.95
$ws.Range("U155").Value2 = 3
$ws.Range("V155").Value2 = 1.8
$ws.Range("W155").Value2 = 2.05
$ws.Range("X155").Value2 = -1
$ws.Range("Y155").Value2 = 4.25
$ws.Range("AA155").Value2 = -1
$ws.Range("AB155").Value2 = 0.95
$ws.Range("AC155").Value2 = 0.8
$ws.Range("AD155").Value2 = -1

# Row 156 <- original row 155
$ws.Range("B156").Value2 = 6810132
$ws.Range("E156").Value2 = "SintTruidense"
$ws.Range("F156").Value2 = "Charleroi"
$ws.Range("G156").Value2 = 1
$ws.Range("H156").Value2 = 0
$ws.Range("I156").Value2 = 0
$ws.Range("K156").Value2 = "H"
$ws.Range("L156").Value2 = 2.2
$ws.Range("M156").Value2 = 3.4
$ws.Range("N156").Value2 = 3.2
$ws.Range("O156").Value2 = 2.3
$ws.Range("P156").Value2 = 3.3
$ws.Range("Q156").Value2 = 3
$ws.Range("R156").Value2 = -0.25
$ws.Range("S156").Value2 = 2
$ws.Range("T156").Value2 = 1.85
$ws.Range("U156").Value2 = 2.25
$ws.Range("V156").Value2 = 1.85
$ws.Range("W156").Value2 = 2
$ws.Range("X156").Value2 = 1.3
$ws.Range("Y156").Value2 = -1
$ws.Range("AA156").Value2 = 1
$ws.Range("AB156").Value2 = -1
$ws.Range("AC156").Value2 = -1
$ws.Range("AD156").Value2 = 1

# Row 159 <- original row 160
$ws.Range("B159").Value2 = 6810145
$ws.Range("E159").Value2 = "KV Kortrijk"
$ws.Range("F159").Value2 = "Gent"
$ws.Range("G159").Value2 = 0
$ws.Range("H159").Value2 = 2
$ws.Range("I159").Value2 = 0
$ws.Range("J159").Value2 = 1
$ws.Range("K159").Value2 = "A"
$ws.Range("L159").Value2 = 7
$ws.Range("M159").Value2 = 5
$ws.Range("N159").Value2 = 1.4
$ws.Range("O159").Value2 = 8.5
$ws.Range("P159").Value2 = 5.75
$ws.Range("Q159").Value2 = 1.285
$ws.Range("R159").Value2 = 1.5
$ws.Range("S159").Value2 = 2.025
$ws.Range("T159").Value2 = 1.825
$ws.Range("U159").Value2 = 3.25
$ws.Range("V159").Value2 = 2.05
$ws.Range("W159").Value2 = 1.8
$ws.Range("X159").Value2 = -1
$ws.Range("Z159").Value2 = 0.2849999999999999
$ws.Range("AA159").Value2 = -1
$ws.Range("AB159").Value2 = 0.825
$ws.Range("AC159").Value2 = -1
$ws.Range("AD159").Value2 = 0.8

# Row 160 <- original row 159
$ws.Range("B160").Value2 = 6810142
$ws.Range("E160").Value2 = "OH Leuven"
$ws.Range("F160").Value2 = "Eupen"
$ws.Range("G160").Value2 = 3
$ws.Range("H160").Value2 = 0
$ws.Range("I160").Value2 = 2
$ws.Range("J160").Value2 = 0
$ws.Range("K160").Value2 = "H"
$ws.Range("L160").Value2 = 1.75
$ws.Range("M160").Value2 = 4
$ws.Range("N160").Value2 = 4
$ws.Range("O160").Value2 = 1.8
$ws.Range("P160").Value2 = 3.8
$ws.Range("Q160").Value2 = 3.8
$ws.Range("R160").Value2 = -0.5
$ws.Range("S160").Value2 = 1.825
$ws.Range("T160").Value2 = 2.025
$ws.Range("U160").Value2 = 3
$ws.Range("V160").Value2 = 1.975
$ws.Range("W160").Value2 = 1.875
$ws.Range("X160").Value2 = 0.8
$ws.Range("Z160").Value2 = -1
$ws.Range("AA160").Value2 = 0.825
$ws.Range("AB160").Value2 = -1
$ws.Range("AC160").Value2 = 0
$ws.Range("AD160").Value2 = 0

# Row 175 <- original row 176
$ws.Range("B175").Value2 = 6810158
$ws.Range("E175").Value2 = "Gent"
$ws.Range("F175").Value2 = "Westerlo"
$ws.Range("G175").Value2 = 2
$ws.Range("H175").Value2 = 2
$ws.Range("I175").Value2 = 1
$ws.Range("J175").Value2 = 1
$ws.Range("K175").Value2 = "D"
$ws.Range("L175").Value2 = 1.363
$ws.Range("M175").Value2 = 5
$ws.Range("N175").Value2 = 6.5
$ws.Range("O175").Value2 = 1.444
$ws.Range("P175").Value2 = 4.5
$ws.Range("Q175").Value2 = 5.75
$ws.Range("R175").Value2 = -1.25
$ws.Range("S175").Value2 = 2
$ws.Range("T175").Value2 = 1.85
$ws.Range("U175").Value2 = 3
$ws.Range("V175").Value2 = 2.025
$ws.Range("W175").Value2 = 1.825
$ws.Range("Y175").Value2 = 3.5
$ws.Range("Z175").Value2 = -1
$ws.Range("AB175").Value2 = 0.8500000000000001
$ws.Range("AC175").Value2 = 1.025

# Row 176 <- original row 175
$ws.Range("B176").Value2 = 6810159
$ws.Range("E176").Value2 = "Charleroi"
$ws.Range("F176").Value2 = "Club Brugge"
$ws.Range("G176").Value2 = 1
$ws.Range("H176").Value2 = 4
$ws.Range("I176").Value2 = 0
$ws.Range("J176").Value2 = 3
$ws.Range("K176").Value2 = "A"
$ws.Range("L176").Value2 = 5.5
$ws.Range("M176").Value2 = 4.2
$ws.Range("N176").Value2 = 1.5
$ws.Range("O176").Value2 = 6
$ws.Range("P176").Value2 = 4.2
$ws.Range("Q176").Value2 = 1.45
$ws.Range("R176").Value2 = 1.25
$ws.Range("S176").Value2 = 1.825
$ws.Range("T176").Value2 = 2.025
$ws.Range("U176").Value2 = 2.75
$ws.Range("V176").Value2 = 1.95
$ws.Range("W176").Value2 = 1.9
$ws.Range("Y176").Value2 = -1
$ws.Range("Z176").Value2 = 0.45
$ws.Range("AB176").Value2 = 1.025
$ws.Range("AC176").Value2 = 0.95

# Row 185 <- original row 186
$ws.Range("B185").Value2 = 6810164
$ws.Range("E185").Value2 = "Union Saint Gilloise"
$ws.Range("F185").Value2 = "RWD Molenbeek"
$ws.Range("G185").Value2 = 3
$ws.Range("H185").Value2 = 2
$ws.Range("I185").Value2 = 2
$ws.Range("K185").Value2 = "H"
$ws.Range("L185").Value2 = 1.2
$ws.Range("M185").Value2 = 7
$ws.Range("N185").Value2 = 12
$ws.Range("O185").Value2 = 1.142
$ws.Range("P185").Value2 = 8.5
$ws.Range("Q185").Value2 = 15
$ws.Range("R185").Value2 = -2.25
$ws.Range("S185").Value2 = 1.925
$ws.Range("T185").Value2 = 1.925
$ws.Range("U185").Value2 = 3.5
$ws.Range("V185").Value2 = 2.025
$ws.Range("W185").Value2 = 1.825
$ws.Range("X185").Value2 = 0.1419999999999999
$ws.Range("Z185").Value2 = -1
$ws.Range("AB185").Value2 = 0.925
$ws.Range("AC185").Value2 = 1.025
$ws.Range("AD185").Value2 = -1

# Row 186 <- original row 185
$ws.Range("B186").Value2 = 6810162
$ws.Range("E186").Value2 = "Standard Liege"
$ws.Range("F186").Value2 = "Antwerp"
$ws.Range("G186").Value2 = 0
$ws.Range("H186").Value2 = 1
$ws.Range("I186").Value2 = 0
$ws.Range("K186").Value2 = "A"
$ws.Range("L186").Value2 = 4
$ws.Range("M186").Value2 = 3.6
$ws.Range("N186").Value2 = 1.85
$ws.Range("O186").Value2 = 3.1
$ws.Range("P186").Value2 = 3.2
$ws.Range("Q186").Value2 = 2.3
$ws.Range("R186").Value2 = 0.25
$ws.Range("S186").Value2 = 1.8
$ws.Range("T186").Value2 = 2.05
$ws.Range("U186").Value2 = 2.25
$ws.Range("V186").Value2 = 1.875
$ws.Range("W186").Value2 = 1.975
$ws.Range("X186").Value2 = -1
$ws.Range("Z186").Value2 = 1.3
$ws.Range("AB186").Value2 = 1.05
$ws.Range("AC186").Value2 = -1
$ws.Range("AD186").Value2 = 0.9750000000000001

# Row 187 <- original row 188
$ws.Range("B187").Value2 = 6810163
$ws.Range("E187").Value2 = "SintTruidense"
$ws.Range("F187").Value2 = "Gent"
$ws.Range("G187").Value2 = 4
$ws.Range("H187").Value2 = 1
$ws.Range("K187").Value2 = "H"
$ws.Range("L187").Value2 = 3.6
$ws.Range("M187").Value2 = 3.6
$ws.Range("N187").Value2 = 1.95
$ws.Range("O187").Value2 = 3.25
$ws.Range("P187").Value2 = 3.4
$ws.Range("Q187").Value2 = 2.15
$ws.Range("S187").Value2 = 1.95
$ws.Range("T187").Value2 = 1.9
$ws.Range("V187").Value2 = 1.975
$ws.Range("W187").Value2 = 1.875
$ws.Range("X187").Value2 = 2.25
$ws.Range("Y187").Value2 = -1
$ws.Range("AA187").Value2 = 0.95
$ws.Range("AB187").Value2 = -1
$ws.Range("AC187").Value2 = 0.9750000000000001

# Row 188 <- original row 187
$ws.Range("B188").Value2 = 6810166
$ws.Range("E188").Value2 = "KV Mechelen"
$ws.Range("F188").Value2 = "Anderlecht"
$ws.Range("G188").Value2 = 2
$ws.Range("H188").Value2 = 2
$ws.Range("K188").Value2 = "D"
$ws.Range("L188").Value2 = 3.5
$ws.Range("M188").Value2 = 3.5
$ws.Range("N188").Value2 = 2
$ws.Range("O188").Value2 = 3
$ws.Range("P188").Value2 = 3.5
$ws.Range("Q188").Value2 = 2.2
$ws.Range("S188").Value2 = 1.925
$ws.Range("T188").Value2 = 1.925
$ws.Range("V188").Value2 = 1.875
$ws.Range("W188").Value2 = 1.975
$ws.Range("X188").Value2 = -1
$ws.Range("Y188").Value2 = 2.5
$ws.Range("AA188").Value2 = 0.4625
$ws.Range("AB188").Value2 = -0.5
$ws.Range("AC188").Value2 = 0.875

# Row 241 <- original row 243
$ws.Range("B241").Value2 = 6810219
$ws.Range("E241").Value2 = "OH Leuven"
$ws.Range("F241").Value2 = "KV Mechelen"
$ws.Range("G241").Value2 = 1
$ws.Range("I241").Value2 = 0
$ws.Range("L241").Value2 = 2.8
$ws.Range("M241").Value2 = 3.5
$ws.Range("N241").Value2 = 2.375
$ws.Range("O241").Value2 = 2.7
$ws.Range("P241").Value2 = 3.5
$ws.Range("Q241").Value2 = 2.45
$ws.Range("R241").Value2 = 0
$ws.Range("S241").Value2 = 2.025
$ws.Range("T241").Value2 = 1.825
$ws.Range("V241").Value2 = 1.925
$ws.Range("W241").Value2 = 1.925
$ws.Range("X241").Value2 = 1.7
$ws.Range("AA241").Value2 = 1.025
$ws.Range("AC241").Value2 = -1
$ws.Range("AD241").Value2 = 0.925

# Row 242 <- original row 241
$ws.Range("B242").Value2 = 6942395
$ws.Range("E242").Value2 = "Gent"
$ws.Range("F242").Value2 = "Charleroi"
$ws.Range("G242").Value2 = 5
$ws.Range("L242").Value2 = 1.571
$ws.Range("M242").Value2 = 4
$ws.Range("N242").Value2 = 5.75
$ws.Range("O242").Value2 = 1.4
$ws.Range("P242").Value2 = 4.333
$ws.Range("Q242").Value2 = 8
$ws.Range("R242").Value2 = -1.25
$ws.Range("S242").Value2 = 2
$ws.Range("T242").Value2 = 1.85
$ws.Range("U242").Value2 = 2.75
$ws.Range("X242").Value2 = 0.3999999999999999
$ws.Range("AA242").Value2 = 1

# Row 243 <- original row 242
$ws.Range("B243").Value2 = 6870199
$ws.Range("E243").Value2 = "Cercle Brugge"
$ws.Range("F243").Value2 = "RWD Molenbeek"
$ws.Range("G243").Value2 = 4
$ws.Range("I243").Value2 = 2
$ws.Range("L243").Value2 = 1.363
$ws.Range("M243").Value2 = 5.5
$ws.Range("N243").Value2 = 7.5
$ws.Range("O243").Value2 = 1.3
$ws.Range("P243").Value2 = 6
$ws.Range("Q243").Value2 = 8.5
$ws.Range("R243").Value2 = -1.75
$ws.Range("U243").Value2 = 3.25
$ws.Range("V243").Value2 = 1.95
$ws.Range("W243").Value2 = 1.9
$ws.Range("X243").Value2 = 0.3
$ws.Range("AC243").Value2 = 0.95
$ws.Range("AD243").Value2 = -1

# Row 275 <- original row 277
$ws.Range("B275").Value2 = 7979470
$ws.Range("E275").Value2 = "Westerlo"
$ws.Range("F275").Value2 = "OH Leuven"
$ws.Range("G275").Value2 = 1
$ws.Range("H275").Value2 = 1
$ws.Range("I275").Value2 = 1
$ws.Range("M275").Value2 = 3.6
$ws.Range("N275").Value2 = 2.6
$ws.Range("O275").Value2 = 2.45
$ws.Range("P275").Value2 = 3.75
$ws.Range("Q275").Value2 = 2.55
$ws.Range("U275").Value2 = 3
$ws.Range("V275").Value2 = 1.85
$ws.Range("W275").Value2 = 2
$ws.Range("Y275").Value2 = 2.75
$ws.Range("AD275").Value2 = 1

# Row 276 <- original row 275
$ws.Range("B276").Value2 = 7979471
$ws.Range("E276").Value2 = "Standard Liege"
$ws.Range("F276").Value2 = "KV Mechelen"
$ws.Range("H276").Value2 = 0
$ws.Range("J276").Value2 = 0
$ws.Range("K276").Value2 = "D"
$ws.Range("L276").Value2 = 2.5
$ws.Range("M276").Value2 = 3.4
$ws.Range("N276").Value2 = 2.7
$ws.Range("O276").Value2 = 2.5
$ws.Range("P276").Value2 = 3.5
$ws.Range("Q276").Value2 = 2.625
$ws.Range("R276").Value2 = 0
$ws.Range("S276").Value2 = 1.875
$ws.Range("T276").Value2 = 1.975
$ws.Range("U276").Value2 = 2.75
$ws.Range("V276").Value2 = 1.925
$ws.Range("W276").Value2 = 1.925
$ws.Range("Y276").Value2 = 2.5
$ws.Range("Z276").Value2 = -1
$ws.Range("AA276").Value2 = 0
$ws.Range("AB276").Value2 = 0
$ws.Range("AD276").Value2 = 0.925

# Row 277 <- original row 276
$ws.Range("B277").Value2 = 7979346
$ws.Range("E277").Value2 = "SintTruidense"
$ws.Range("F277").Value2 = "Gent"
$ws.Range("G277").Value2 = 0
$ws.Range("H277").Value2 = 2
$ws.Range("I277").Value2 = 0
$ws.Range("J277").Value2 = 1
$ws.Range("K277").Value2 = "A"
$ws.Range("L277").Value2 = 3.6
$ws.Range("M277").Value2 = 3.5
$ws.Range("N277").Value2 = 2
$ws.Range("O277").Value2 = 3.3
$ws.Range("P277").Value2 = 3.6
$ws.Range("Q277").Value2 = 2.05
$ws.Range("R277").Value2 = 0.25
$ws.Range("S277").Value2 = 2.025
$ws.Range("T277").Value2 = 1.825
$ws.Range("V277").Value2 = 1.975
$ws.Range("W277").Value2 = 1.875
$ws.Range("Y277").Value2 = -1
$ws.Range("Z277").Value2 = 1.05
$ws.Range("AA277").Value2 = -1
$ws.Range("AB277").Value2 = 0.825
$ws.Range("AD277").Value2 = 0.875

# Row 278 <- original row 279
$ws.Range("B278").Value2 = 7979357
$ws.Range("E278").Value2 = "Club Brugge"
$ws.Range("F278").Value2 = "Genk"
$ws.Range("G278").Value2 = 4
$ws.Range("I278").Value2 = 1
$ws.Range("L278").Value2 = 1.85
$ws.Range("M278").Value2 = 3.75
$ws.Range("N278").Value2 = 3.9
$ws.Range("O278").Value2 = 1.75
$ws.Range("P278").Value2 = 3.75
$ws.Range("Q278").Value2 = 4.5
$ws.Range("R278").Value2 = -0.75
$ws.Range("S278").Value2 = 2
$ws.Range("T278").Value2 = 1.85
$ws.Range("V278").Value2 = 2.025
$ws.Range("W278").Value2 = 1.825
$ws.Range("X278").Value2 = 0.75
$ws.Range("AA278").Value2 = 1
$ws.Range("AC278").Value2 = 1.025
$ws.Range("AD278").Value2 = -1

# Row 279 <- original row 278
$ws.Range("B279").Value2 = 7979473
$ws.Range("E279").Value2 = "Anderlecht"
$ws.Range("F279").Value2 = "Cercle Brugge"
$ws.Range("G279").Value2 = 3
$ws.Range("I279").Value2 = 2
$ws.Range("L279").Value2 = 1.909
$ws.Range("M279").Value2 = 3.6
$ws.Range("N279").Value2 = 3.8
$ws.Range("O279").Value2 = 1.8
$ws.Range("P279").Value2 = 3.8
$ws.Range("Q279").Value2 = 4
$ws.Range("R279").Value2 = -0.5
$ws.Range("S279").Value2 = 1.85
$ws.Range("T279").Value2 = 2
$ws.Range("V279").Value2 = 1.85
$ws.Range("W279").Value2 = 2
$ws.Range("X279").Value2 = 0.8
$ws.Range("AA279").Value2 = 0.8500000000000001
$ws.Range("AC279").Value2 = 0.425
$ws.Range("AD279").Value2 = -0.5

# Row 310 <- original row 311
$ws.Range("B310").Value2 = 8009913
$ws.Range("E310").Value2 = "OH Leuven"
$ws.Range("F310").Value2 = "Westerlo"
$ws.Range("G310").Value2 = 1
$ws.Range("I310").Value2 = 0
$ws.Range("J310").Value2 = 2
$ws.Range("K310").Value2 = "A"
$ws.Range("L310").Value2 = 2.05
$ws.Range("M310").Value2 = 3.7
$ws.Range("N310").Value2 = 3.2
$ws.Range("O310").Value2 = 1.909
$ws.Range("P310").Value2 = 3.8
$ws.Range("Q310").Value2 = 3.4
$ws.Range("R310").Value2 = -0.5
$ws.Range("S310").Value2 = 1.975
$ws.Range("T310").Value2 = 1.875
$ws.Range("U310").Value2 = 3.25
$ws.Range("V310").Value2 = 1.875
$ws.Range("W310").Value2 = 1.975
$ws.Range("X310").Value2 = -1
$ws.Range("Z310").Value2 = 2.4
$ws.Range("AA310").Value2 = -1
$ws.Range("AB310").Value2 = 0.875
$ws.Range("AC310").Value2 = -0.5
$ws.Range("AD310").Value2 = 0.4875

# Row 311 <- original row 312
$ws.Range("B311").Value2 = 8009351
$ws.Range("E311").Value2 = "Gent"
$ws.Range("F311").Value2 = "SintTruidense"
$ws.Range("G311").Value2 = 2
$ws.Range("H311").Value2 = 0
$ws.Range("J311").Value2 = 0
$ws.Range("K311").Value2 = "H"
$ws.Range("L311").Value2 = 1.5
$ws.Range("M311").Value2 = 4.5
$ws.Range("N311").Value2 = 5.25
$ws.Range("O311").Value2 = 1.533
$ws.Range("P311").Value2 = 4.5
$ws.Range("Q311").Value2 = 5
$ws.Range("R311").Value2 = -1
$ws.Range("S311").Value2 = 1.9
$ws.Range("T311").Value2 = 1.95
$ws.Range("U311").Value2 = 3.5
$ws.Range("V311").Value2 = 1.95
$ws.Range("W311").Value2 = 1.9
$ws.Range("X311").Value2 = 0.5329999999999999
$ws.Range("Z311").Value2 = -1
$ws.Range("AA311").Value2 = 0.8999999999999999
$ws.Range("AB311").Value2 = -1
$ws.Range("AC311").Value2 = -1
$ws.Range("AD311").Value2 = 0.8999999999999999

# Row 312 <- original row 310
$ws.Range("B312").Value2 = 8009914
$ws.Range("E312").Value2 = "KV Mechelen"
$ws.Range("F312").Value2 = "Standard Liege"
$ws.Range("G312").Value2 = 3
$ws.Range("H312").Value2 = 2
$ws.Range("I312").Value2 = 2
$ws.Range("L312").Value2 = 1.909
$ws.Range("M312").Value2 = 3.6
$ws.Range("N312").Value2 = 3.6
$ws.Range("O312").Value2 = 1.65
$ws.Range("P312").Value2 = 4.333
$ws.Range("Q312").Value2 = 4.333
$ws.Range("R312").Value2 = -0.75
$ws.Range("S312").Value2 = 1.825
$ws.Range("T312").Value2 = 2.025
$ws.Range("V312").Value2 = 1.925
$ws.Range("W312").Value2 = 1.925
$ws.Range("X312").Value2 = 0.6499999999999999
$ws.Range("AA312").Value2 = 0.4125
$ws.Range("AB312").Value2 = -0.5
$ws.Range("AC312").Value2 = 0.925
$ws.Range("AD312").Value2 = -1

# Row 313 <- original row 314
$ws.Range("B313").Value2 = 8009865
$ws.Range("E313").Value2 = "Union Saint Gilloise"
$ws.Range("F313").Value2 = "Genk"
$ws.Range("G313").Value2 = 2
$ws.Range("K313").Value2 = "H"
$ws.Range("L313").Value2 = 1.666
$ws.Range("M313").Value2 = 3.75
$ws.Range("N313").Value2 = 4.333
$ws.Range("O313").Value2 = 1.5
$ws.Range("P313").Value2 = 4.2
$ws.Range("Q313").Value2 = 5
$ws.Range("R313").Value2 = -1
$ws.Range("S313").Value2 = 1.875
$ws.Range("T313").Value2 = 1.975
$ws.Range("V313").Value2 = 1.875
$ws.Range("W313").Value2 = 1.975
$ws.Range("X313").Value2 = 0.5
$ws.Range("Y313").Value2 = -1
$ws.Range("AA313").Value2 = 0.875
$ws.Range("AB313").Value2 = -1
$ws.Range("AD313").Value2 = 0.9750000000000001

# Row 314 <- original row 313
$ws.Range("B314").Value2 = 8009325
$ws.Range("E314").Value2 = "Club Brugge"
$ws.Range("F314").Value2 = "Cercle Brugge"
$ws.Range("G314").Value2 = 0
$ws.Range("K314").Value2 = "D"
$ws.Range("L314").Value2 = 1.444
$ws.Range("M314").Value2 = 4.5
$ws.Range("N314").Value2 = 5.5
$ws.Range("O314").Value2 = 1.615
$ws.Range("P314").Value2 = 3.8
$ws.Range("Q314").Value2 = 4.5
$ws.Range("R314").Value2 = -0.75
$ws.Range("S314").Value2 = 1.825
$ws.Range("T314").Value2 = 2.025
$ws.Range("V314").Value2 = 2.05
$ws.Range("W314").Value2 = 1.8
$ws.Range("X314").Value2 = -1
$ws.Range("Y314").Value2 = 2.8
$ws.Range("AA314").Value2 = -1
$ws.Range("AB314").Value2 = 1.025
$ws.Range("AD314").Value2 = 0.8
